$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.762.11"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "3.258.51"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.257.76"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.506"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "3.799.02"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("D16").Value = "66.859.33"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "3.245.65"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.752"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("B27").Value = "Hedera"
$ws.Range("C27").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.140"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +59.10%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +23.20%  "
$ws.Range("D38").Value = "0.0₃0798"
$ws.Range("E38").Value = "  +18.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "495.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +4.34%  "
$ws.Range("D46").Value = "2.978.36"
$ws.Range("E46").Value = "  +5.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.03%  "
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
